# Applies the "Bug fix + added new heuristic + experiments" edit described
# by the target diff:
#  - experiment_description: new experiment 6 ("Transition driven: 1 +
#    hospital_less_than_14_days heuristic") inserted before experiment 1001
#  - experiment_specification: matching rows for experiment 6 inserted
#    before the 1001 block
#  - run_description: two new runs (5, 6) appended
#  - run_specification: run 5 and run 6 experiment mappings appended
#  - heuristics_description: new heuristic 5 (hospital_less_than_14_days)
#    appended

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# run_description: append run 5 and run 6 (their text is referenced by
# the experiment_description sheet further down, so create them first to
# mirror the author's original shared-string ordering).
# ---------------------------------------------------------------------
$wsRunDesc = $wb.Worksheets.Item("run_description")
$wsRunDesc.Cells.Item(6, 1).Value = 5
$wsRunDesc.Cells.Item(6, 2).Value = "Comparing Base and CA"

# ---------------------------------------------------------------------
# heuristics_description: append heuristic 5.
# ---------------------------------------------------------------------
$wsHeur = $wb.Worksheets.Item("heuristics_description")
$wsHeur.Cells.Item(6, 1).Value = 5
$wsHeur.Cells.Item(6, 2).Value = "hospital_less_than_14_days: If home and worst state is home and next state is not recovered length of stay has to be less than 14 days."

# ---------------------------------------------------------------------
# experiment_description: insert experiment 6 as a new row before the
# existing row for experiment 1001 (old row 7 -> new row 8).
# ---------------------------------------------------------------------
$wsExpDesc = $wb.Worksheets.Item("experiment_description")
$wsExpDesc.Rows.Item(7).Insert()
$wsExpDesc.Cells.Item(7, 1).Value = 6
$wsExpDesc.Cells.Item(7, 2).Value = "Transition driven: 1 + hospital_less_than_14_days heuristic"
$wsExpDesc.Cells.Item(7, 3).Value = "base"
$wsExpDesc.Cells.Item(7, 4).Value = "1;2;3;5"

# ---------------------------------------------------------------------
# run_description: append run 6, referencing the new experiment.
# ---------------------------------------------------------------------
$wsRunDesc.Cells.Item(7, 1).Value = 6
$wsRunDesc.Cells.Item(7, 2).Value = "Comparing aggressive transition and los driven base models"

# ---------------------------------------------------------------------
# experiment_specification: insert the 3 rows (home / inpatient_ward /
# intensive_care_unit) describing experiment 6, right before the 1001
# block (old row 17 -> new row 20).
# ---------------------------------------------------------------------
$wsExpSpec = $wb.Worksheets.Item("experiment_specification")
$wsExpSpec.Rows.Item(17).Insert()
$wsExpSpec.Rows.Item(17).Insert()
$wsExpSpec.Rows.Item(17).Insert()

$newSpecRows = @(
    @(6, "home", "none", "age_simple", "age_simple"),
    @(6, "inpatient_ward", "none", "age_simple", "none"),
    @(6, "intensive_care_unit", "none", "age_simple", "none")
)
for ($i = 0; $i -lt $newSpecRows.Count; $i++) {
    $r = 17 + $i
    $rowValues = $newSpecRows[$i]
    for ($j = 0; $j -lt $rowValues.Count; $j++) {
        $wsExpSpec.Cells.Item($r, $j + 1).Value = $rowValues[$j]
    }
}

# ---------------------------------------------------------------------
# run_specification: append the experiment mappings for run 5 and run 6.
# ---------------------------------------------------------------------
$wsRunSpec = $wb.Worksheets.Item("run_specification")
$newRunSpecRows = @(
    @(5, 1),
    @(5, 1001),
    @(5, 4),
    @(5, 1003),
    @(6, 5),
    @(6, 6)
)
for ($i = 0; $i -lt $newRunSpecRows.Count; $i++) {
    $r = 12 + $i
    $rowValues = $newRunSpecRows[$i]
    $wsRunSpec.Cells.Item($r, 1).Value = $rowValues[0]
    $wsRunSpec.Cells.Item($r, 2).Value = $rowValues[1]
}

# ---------------------------------------------------------------------
# Selection / active-sheet bookkeeping to mirror the end-user state
# captured in the diff (run_specification ends up the active sheet).
# ---------------------------------------------------------------------
$wsExpDesc.Activate()
$wsExpDesc.Range("B20").Select()

$wsExpSpec.Activate()
$wsExpSpec.Range("D21").Select()

$wsRunDesc.Activate()
$wsRunDesc.Range("B25").Select()

$wsHeur.Activate()
$wsHeur.Range("B11").Select()

$wsRunSpec.Activate()
$wsRunSpec.Range("C22").Select()
